$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (G1:I1) matching the style of the existing
# bold header row, plus the new data value in G2.
$ws.Range("G1").Value = "Hz"
$ws.Range("H1").Value = "V Thresh [V] @ .5s"
$ws.Range("I1").Value = "W thresh [s] @ 1.5 Vt"

$ws.Range("G2").Value = 4

# Apply the same bold header formatting used by the existing header cells
$ws.Range("G1:I1").Font.Bold = $true

# Update the selection to match the newly added range
$ws.Range("G1:I2").Select()
